$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2876.7778
$ws.Cells.Item(19, 9).Value = 2478.8
$ws.Cells.Item(19, 11).Value = 2478.8
$ws.Cells.Item(19, 13).Value = -2303.8

$ws.Cells.Item(41, 8).Value = 583.4286
$ws.Cells.Item(41, 9).Value = 480.66666
$ws.Cells.Item(41, 11).Value = 480.66666
$ws.Cells.Item(41, 13).Value = -40.66665999999998

$ws.Cells.Item(86, 8).Value = 5712.091
$ws.Cells.Item(86, 9).Value = 5547.8335
$ws.Cells.Item(86, 10).Value = 5909.2
$ws.Cells.Item(86, 11).Value = 5547.8335
$ws.Cells.Item(86, 12).Value = 5909.2
$ws.Cells.Item(86, 13).Value = -4424.8335
$ws.Cells.Item(86, 14).Value = -8155.2

$ws.Cells.Item(89, 8).Value = 5712.091
$ws.Cells.Item(89, 9).Value = 5547.8335
$ws.Cells.Item(89, 10).Value = 5909.2
$ws.Cells.Item(89, 11).Value = 27739.1675
$ws.Cells.Item(89, 12).Value = 29546
$ws.Cells.Item(89, 13).Value = -22123.1675
$ws.Cells.Item(89, 14).Value = -40778

$ws.Cells.Item(125, 8).Value = 3349
$ws.Cells.Item(125, 9).Value = 1903.3334
$ws.Cells.Item(125, 10).Value = 5517.5
$ws.Cells.Item(125, 11).Value = 17130.0006
$ws.Cells.Item(125, 12).Value = 49657.5
$ws.Cells.Item(125, 13).Value = -14670.0006
$ws.Cells.Item(125, 14).Value = -54577.5

$ws.Cells.Item(132, 8).Value = 2321
$ws.Cells.Item(132, 9).Value = 1996.4445
$ws.Cells.Item(132, 11).Value = 5989.333500000001
$ws.Cells.Item(132, 13).Value = -3459.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2474
$ws.Cells.Item(2, 9).Value = 1459.8
$ws.Cells.Item(2, 10).Value = 3741.75
$ws.Cells.Item(2, 11).Value = 1459.8
$ws.Cells.Item(2, 12).Value = 3741.75
$ws.Cells.Item(2, 13).Value = -1346.8
$ws.Cells.Item(2, 14).Value = -3967.75

$ws.Cells.Item(28, 8).Value = 10262.714
$ws.Cells.Item(28, 9).Value = 10262.714
$ws.Cells.Item(28, 11).Value = 10262.714
$ws.Cells.Item(28, 13).Value = -10070.714

$ws.Cells.Item(74, 8).Value = 804.4286
$ws.Cells.Item(74, 9).Value = 804.4286
$ws.Cells.Item(74, 11).Value = 804.4286
$ws.Cells.Item(74, 13).Value = 69.57140000000004

$ws.Cells.Item(77, 8).Value = 804.4286
$ws.Cells.Item(77, 9).Value = 804.4286
$ws.Cells.Item(77, 11).Value = 4022.143
$ws.Cells.Item(77, 13).Value = 345.857

$ws.Cells.Item(99, 8).Value = 10262.714
$ws.Cells.Item(99, 9).Value = 10262.714
$ws.Cells.Item(99, 11).Value = 10262.714
$ws.Cells.Item(99, 13).Value = -7267.714

$ws.Cells.Item(116, 8).Value = 2474
$ws.Cells.Item(116, 9).Value = 1459.8
$ws.Cells.Item(116, 10).Value = 3741.75
$ws.Cells.Item(116, 11).Value = 1459.8
$ws.Cells.Item(116, 12).Value = 3741.75
$ws.Cells.Item(116, 13).Value = 834.2
$ws.Cells.Item(116, 14).Value = -8329.75

$ws.Cells.Item(132, 8).Value = 1895
$ws.Cells.Item(132, 9).Value = 1895
$ws.Cells.Item(132, 11).Value = 5685
$ws.Cells.Item(132, 13).Value = -3155

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2474
$ws.Cells.Item(3, 9).Value = 1459.8
$ws.Cells.Item(3, 10).Value = 3741.75
$ws.Cells.Item(3, 11).Value = 1459.8
$ws.Cells.Item(3, 12).Value = 3741.75
$ws.Cells.Item(3, 13).Value = -1345.8
$ws.Cells.Item(3, 14).Value = -3969.75

$ws.Cells.Item(36, 8).Value = 2425
$ws.Cells.Item(36, 9).Value = 2425
$ws.Cells.Item(36, 11).Value = 2425
$ws.Cells.Item(36, 13).Value = -1891

$ws.Cells.Item(74, 8).Value = 52922.668
$ws.Cells.Item(74, 10).Value = 52922.668
$ws.Cells.Item(74, 12).Value = 52922.668
$ws.Cells.Item(74, 14).Value = -54794.668

$ws.Cells.Item(77, 8).Value = 52922.668
$ws.Cells.Item(77, 10).Value = 52922.668
$ws.Cells.Item(77, 12).Value = 158768.004
$ws.Cells.Item(77, 14).Value = -168128.004

$ws.Cells.Item(96, 8).Value = 8571
$ws.Cells.Item(96, 9).Value = 8571
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 8571
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -5825
$ws.Cells.Item(96, 14).ClearContents()

$ws.Cells.Item(105, 8).Value = 3129.5
$ws.Cells.Item(105, 9).Value = 3129.5
$ws.Cells.Item(105, 11).Value = 3129.5
$ws.Cells.Item(105, 13).Value = -1382.5

$ws.Cells.Item(109, 8).Value = 90796
$ws.Cells.Item(109, 10).Value = 90796
$ws.Cells.Item(109, 12).Value = 90796
$ws.Cells.Item(109, 14).Value = -93570

$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2087.5
$ws.Cells.Item(31, 9).Value = 1116.6666
$ws.Cells.Item(31, 10).Value = 5000
$ws.Cells.Item(31, 11).Value = 1116.6666
$ws.Cells.Item(31, 12).Value = 5000
$ws.Cells.Item(31, 13).Value = -821.6666
$ws.Cells.Item(31, 14).Value = -5590

$ws.Cells.Item(34, 8).Value = 2087.5
$ws.Cells.Item(34, 9).Value = 1116.6666
$ws.Cells.Item(34, 10).Value = 5000
$ws.Cells.Item(34, 11).Value = 1116.6666
$ws.Cells.Item(34, 12).Value = 5000
$ws.Cells.Item(34, 13).Value = -914.6666
$ws.Cells.Item(34, 14).Value = -5404

$ws.Cells.Item(107, 8).Value = 396.125
$ws.Cells.Item(107, 9).Value = 358.92856
$ws.Cells.Item(107, 10).Value = 656.5
$ws.Cells.Item(107, 11).Value = 358.92856
$ws.Cells.Item(107, 12).Value = 656.5
$ws.Cells.Item(107, 13).Value = 1561.07144
$ws.Cells.Item(107, 14).Value = -4496.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1036.3334
$ws.Cells.Item(107, 9).Value = 700
$ws.Cells.Item(107, 10).Value = 1204.5
$ws.Cells.Item(107, 11).Value = 2100
$ws.Cells.Item(107, 12).Value = 3613.5
$ws.Cells.Item(107, 13).Value = -180
$ws.Cells.Item(107, 14).Value = -7453.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 151.33333
$ws.Cells.Item(13, 9).Value = 247
$ws.Cells.Item(13, 10).Value = 103.5
$ws.Cells.Item(13, 11).Value = 247
$ws.Cells.Item(13, 12).Value = 103.5
$ws.Cells.Item(13, 13).Value = -108
$ws.Cells.Item(13, 14).Value = -381.5

$ws.Cells.Item(98, 8).Value = 21252.5
$ws.Cells.Item(98, 10).Value = 21252.5
$ws.Cells.Item(98, 12).Value = 21252.5
$ws.Cells.Item(98, 14).Value = -27242.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1027.2858
$ws.Cells.Item(16, 9).Value = 1027.2858
$ws.Cells.Item(16, 11).Value = 1027.2858
$ws.Cells.Item(16, 13).Value = -857.2858000000001

$ws.Cells.Item(35, 8).Value = 1205.3334
$ws.Cells.Item(35, 9).Value = 1205.3334
$ws.Cells.Item(35, 11).Value = 1205.3334
$ws.Cells.Item(35, 13).Value = -869.3334

$ws.Cells.Item(61, 8).Value = 1622.3334
$ws.Cells.Item(61, 9).Value = 1378.2
$ws.Cells.Item(61, 11).Value = 1378.2
$ws.Cells.Item(61, 13).Value = -1176.2

$ws.Cells.Item(99, 8).Value = 20665
$ws.Cells.Item(99, 9).Value = 20665
$ws.Cells.Item(99, 11).Value = 20665
$ws.Cells.Item(99, 13).Value = -17670

$ws.Cells.Item(113, 8).Value = 1622.3334
$ws.Cells.Item(113, 9).Value = 1378.2
$ws.Cells.Item(113, 11).Value = 1378.2
$ws.Cells.Item(113, 13).Value = 791.8

$ws.Cells.Item(136, 8).Value = 2002.762
$ws.Cells.Item(136, 9).Value = 2003.3158
$ws.Cells.Item(136, 10).Value = 1997.5
$ws.Cells.Item(136, 11).Value = 6009.9474
$ws.Cells.Item(136, 12).Value = 5992.5
$ws.Cells.Item(136, 13).Value = -3459.9474
$ws.Cells.Item(136, 14).Value = -11092.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2649
$ws.Cells.Item(122, 9).Value = 2284.7144
$ws.Cells.Item(122, 10).Value = 3499
$ws.Cells.Item(122, 11).Value = 6854.1432
$ws.Cells.Item(122, 12).Value = 10497
$ws.Cells.Item(122, 13).Value = -4404.1432
$ws.Cells.Item(122, 14).Value = -15397

$ws.Cells.Item(132, 8).Value = 2412
$ws.Cells.Item(132, 9).Value = 2029.7142
$ws.Cells.Item(132, 10).Value = 3750
$ws.Cells.Item(132, 11).Value = 6089.142599999999
$ws.Cells.Item(132, 12).Value = 11250
$ws.Cells.Item(132, 13).Value = -3559.142599999999
$ws.Cells.Item(132, 14).Value = -16310
